# Layout and tooltips update.
$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("categories")

# Update tooltip text for "All Other" category
$ws2.Range("C8").Value = "Including education costs, insurance, social care, financial services"

# Update tooltip text for Communication and electronics
$ws2.Range("C6").Value = "Including mobile phone and internet services, and related electronic devices"

# Update category name for "All Other"
$ws2.Range("A8").Value = "All Other"

# Update tooltip text for Savings
$ws2.Range("C9").Value = "Percentage of income to save (after pension contribution deducted)"

# Update tooltip text for Pension
$ws2.Range("C10").Value = "Percentage of income to save for a pension"

# Make "categories" sheet the active tab and select A16
$ws2.Activate()
$ws2.Range("A16").Select()
